$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('AF7').Value = 'response component cell ontology ID'
$ws.Range('AG7').Value = 'response component marker protein ontology ID'
$ws.Range('AJ7').Value = 'ID of observation within a publication (PMID) and for its submission type '
$ws.Range('AK7').Value = 'ID of observation within its submission type'
$ws.Range('J7').Value = 'tissue cell ontology ID'
$ws.Range('K7').Value = 'response component (original curated cell type)'
$ws.Range('AF8').Value = 'CL:0000946'
$ws.Range('AF9').Value = 'CL:0000946'
$ws.Range('AF10').Value = 'CL:0000946'
$ws.Range('AF11').Value = 'CL:0000946'
$ws.Range('AF12').Value = 'CL:0000946'
$ws.Range('AF13').Value = 'CL:0000904'
$ws.Range('AF14').Value = 'CL:0000904'
$ws.Range('AF15').Value = 'CL:0000900'
$ws.Range('AF16').Value = 'CL:0000904'
$ws.Range('AF17').Value = 'CL:0000904'
$ws.Range('AF18').Value = 'CL:0000900'
$ws.Range('AF19').Value = 'CL:0000980'
$ws.Range('AF20').Value = 'CL:0000980'
$ws.Range('AF21').Value = 'CL:0000980'
$ws.Range('AF22').Value = 'CL:0000623'
$ws.Range('AF23').Value = 'CL:0000979'
$ws.Range('AG23').Value = 'PR:000001020, PR:000001002, PR:000001963'
$ws.Range('AF24').Value = 'CL:0000782'
$ws.Range('AG24').Value = 'PR:000001412'
$ws.Range('AF25').Value = 'CL:0000784'
$ws.Range('AG25').Value = 'PR:P29965, PR:000001412 '
$ws.Range('AF26').Value = 'CL:0000084'
$ws.Range('AF27').Value = 'CL:0000084'
$ws.Range('AF28').Value = 'CL:0000905'
$ws.Range('AF29').Value = 'CL:0000236'
$ws.Range('AF30').Value = 'CL:0002057'
$ws.Range('AF31').Value = 'CL:0002057'
$ws.Range('AF32').Value = 'CL:0000236'
$ws.Range('AF33').Value = 'CL:0000792'
$ws.Range('AF34').Value = 'CL:0000786'
$ws.Range('AF35').Value = 'CL:0000084'
$ws.Range('AF36').Value = 'CL:0000084'
$ws.Range('AF37').Value = 'CL:0000084'
$ws.Range('AF38').Value = 'CL:0000084'
$ws.Range('AF39').Value = 'CL:0000624'
$ws.Range('AF40').Value = 'CL:0000624'
$ws.Range('AF41').Value = 'CL:0000624'
$ws.Range('AF42').Value = 'CL:0000624'
$ws.Range('AF43').Value = 'CL:0000625'
$ws.Range('AF44').Value = 'CL:0000625'
$ws.Range('AF45').Value = 'CL:0000738'
$ws.Range('AF46').Value = 'CL:0000979'
$ws.Range('AG46').Value = 'PR:000001020, PR:000001002, PR:000001963'
$ws.Range('AF47').Value = 'CL:0000624'
$ws.Range('AG47').Value = 'PR:Q96D21 '
$ws.Range('AF48').Value = 'CL:0000787'
$ws.Range('AF49').Value = 'CL:0000982'
$ws.Range('AF50').Value = 'CL:0000576'
$ws.Range('AF51').Value = 'CL:0000980'
$ws.Range('AF52').Value = 'CL:0000980'
$ws.Range('AF53').Value = 'CL:0000980'
$ws.Range('AF54').Value = 'CL:0000624'
$ws.Range('AG54').Value = 'PR:000001383, PR:000025670 '
$ws.Range('AF55').Value = 'CL:0000897'
$ws.Range('AG55').Value = 'PR:000001963, PR:000001203 '
$ws.Range('AF56').Value = 'CL:0000982'
$ws.Range('AF57').Value = 'CL:0000983'
$ws.Range('AF58').Value = 'CL:0000984'
$ws.Range('AF59').Value = 'CL:0000982'
$ws.Range('AF60').Value = 'CL:0000982'
$ws.Range('AF61').Value = 'CL:0000788'
$ws.Range('AG61').Value = 'PR:000001002, PR:000001963'
$ws.Range('AF62').Value = 'CL:0000980'
$ws.Range('AG62').Value = 'PR:000001002, PR:000001963, PR:000001408'
$ws.Range('AF63').Value = 'CL:0000980'
$ws.Range('AF64').Value = 'CL:0000980'
$ws.Range('AF65').Value = 'CL:0000786'
$ws.Range('AF66').Value = 'CL:0000980'
$ws.Range('AF67').Value = 'CL:0000623'
$ws.Range('AF68').Value = 'CL:0000939'
$ws.Range('AG68').Value = 'PR:000001483'
$ws.Range('AF69').Value = 'CL:0000623'
$ws.Range('AF70').Value = 'CL:0000623'
$ws.Range('AF71').Value = 'CL:0000623'
$ws.Range('AF72').Value = 'CL:0000623'
$ws.Range('AF73').Value = 'CL:0000623'
$ws.Range('AF74').Value = 'CL:0000623'
$ws.Range('AF75').Value = 'CL:0000623'
$ws.Range('AF76').Value = 'CL:0000623'
$ws.Range('AF77').Value = 'CL:0000623'
$ws.Range('AF78').Value = 'CL:0000623'
$ws.Range('AF79').Value = 'CL:0000939'
$ws.Range('AG79').Value = 'PR:000001483'
$ws.Range('AF80').Value = 'CL:0000623'
$ws.Range('AF81').Value = 'CL:0000576'
$ws.Range('AF82').Value = 'CL:0000576'
$ws.Range('AF83').Value = 'CL:0000576'
$ws.Range('AF84').Value = 'CL:0000576'
$ws.Range('AF85').Value = 'CL:0000576'
$ws.Range('AF86').Value = 'CL:0000576'
$ws.Range('AF87').Value = 'CL:0000576'
$ws.Range('AF88').Value = 'CL:0000576'
$ws.Range('AF89').Value = 'CL:0000576'
$ws.Range('AF90').Value = 'CL:0002057'
$ws.Range('AF91').Value = 'CL:0002397'
$ws.Range('AF92').Value = 'CL:0000623'
$ws.Range('AG92').Value = 'PR:000001024, PR:000001483'
$ws.Range('AF93').Value = 'CL:0002396'
$ws.Range('AG93').Value = 'PR:000001889, PR:000001483'
$ws.Range('AF94').Value = 'CL:0002057'
$ws.Range('AF95').Value = 'CL:0002397'
$ws.Range('AF96').Value = 'CL:0000896'
$ws.Range('AF97').Value = 'CL:0000775'
$ws.Range('AF98').Value = 'CL:0000624'
$ws.Range('AF99').Value = 'CL:0000624'
$ws.Range('AF100').Value = 'CL:0000624'
$ws.Range('AF101').Value = 'CL:0000738'
$ws.Range('AF102').Value = 'CL:0000542'
$ws.Range('AF103').Value = 'CL:0000784'
$ws.Range('AF104').Value = 'CL:0000623'
$ws.Range('AF105').Value = 'CL:0000236'
$ws.Range('AG105').Value = 'PR:000001020, PR:000001002, PR:000001289, PR:000001963'
$ws.Range('AF106').Value = 'CL:0001054'
$ws.Range('AF107').Value = 'CL:0000624'
$ws.Range('AF108').Value = 'CL:0000784'
$ws.Range('AF109').Value = 'CL:0000623'
$ws.Range('AF110').Value = 'CL:0000624'
$ws.Range('AF111').Value = 'CL:0000451'
$ws.Range('AF112').Value = 'CL:0000236'
$ws.Range('AG112').Value = 'PR:000001020, PR:000001002, PR:000001289'
$ws.Range('AF113').Value = 'CL:0000623'
$ws.Range('AF114').Value = 'CL:0000624'
$ws.Range('AF115').Value = 'CL:0000451'
$ws.Range('AF116').Value = 'CL:0000236'
$ws.Range('AG116').Value = 'PR:000001020, PR:000001002, PR:000001289'
$ws.Range('AF117').Value = 'CL:0000451'
$ws.Range('AF118').Value = 'CL:0000451'
$ws.Range('AF119').Value = 'CL:0000623'
$ws.Range('AF120').Value = 'CL:0000623'
$ws.Range('AF121').Value = 'CL:0000625'
$ws.Range('AF122').Value = 'CL:0000625'
$ws.Range('AF123').Value = 'CL:0000782'
$ws.Range('AG123').Value = 'PR:000001412'
$ws.Range('AF124').Value = 'CL:0000782'
$ws.Range('AG124').Value = 'PR:000001412'
$ws.Range('AF125').Value = 'CL:0000782'
$ws.Range('AG125').Value = 'PR:000001412'
$ws.Range('AF126').Value = 'CL:0000782'
$ws.Range('AG126').Value = 'PR:000001412'
$ws.Range('AF127').Value = 'CL:0000784'
$ws.Range('AG127').Value = 'PR:000001412 '
$ws.Range('AF128').Value = 'CL:0000624'
$ws.Range('AF129').Value = 'CL:0000623'
$ws.Range('AF130').Value = 'CL:0000624'
$ws.Range('AF131').Value = 'CL:0000624'
$ws.Range('AF132').Value = 'CL:0000625'
$ws.Range('AF133').Value = 'CL:0000624'
$ws.Range('AF134').Value = 'CL:0000625'
$ws.Range('AF135').Value = 'CL:0000815'
$ws.Range('AG135').Value = 'PR:000001350 '
$ws.Range('AF136').Value = 'CL:0000815'
$ws.Range('AG136').Value = 'PR:000001350 '
$ws.Range('AF137').Value = 'CL:0000623'
$ws.Range('AF138').Value = 'CL:0000897'
$ws.Range('AF139').Value = 'CL:0000897'
$ws.Range('AF140').Value = 'CL:0000897'
$ws.Range('AF141').Value = 'CL:0000897'
$ws.Range('AF142').Value = 'CL:0000897'
$ws.Range('AF143').Value = 'CL:0000897'
$ws.Range('AF144').Value = 'CL:0000897'
$ws.Range('AF145').Value = 'CL:0000897'
$ws.Range('AG145').Value = 'PR:000001379, PR:000002307 '
$ws.Range('AF146').Value = 'CL:0000897'
$ws.Range('AG146').Value = 'PR:000001379, PR:000002307 '
$ws.Range('AF147').Value = 'CL:0000897'
$ws.Range('AG147').Value = 'PR:000001379, PR:000002307 '
$ws.Range('AF148').Value = 'CL:0000897'
$ws.Range('AG148').Value = 'PR:000001379, PR:000002307 '
$ws.Range('AF149').Value = 'CL:0000625'
$ws.Range('AF150').Value = 'CL:0000625'
$ws.Range('AF151').Value = 'CL:0000625'
$ws.Range('AF152').Value = 'CL:0000625'
$ws.Range('AF153').Value = 'CL:0000625'
$ws.Range('AF154').Value = 'CL:0000624'
$ws.Range('AF155').Value = 'CL:0000794'
$ws.Range('AF156').Value = 'CL:0000775'
$ws.Range('AF157').Value = 'CL:0000576'
$ws.Range('AF158').Value = 'CL:0002397'
$ws.Range('AF159').Value = 'CL:0000623'
$ws.Range('AG159').Value = 'PR:000001024, PR:000001483'
$ws.Range('AF160').Value = 'CL:0000236'
$ws.Range('AG160').Value = 'PR:000001002, PR:000001343'
$ws.Range('AF161').Value = 'CL:0002397'
$ws.Range('AF162').Value = 'CL:0002397'
$ws.Range('AF163').Value = 'CL:0000624'
$ws.Range('AF164').Value = 'CL:0001054'
$ws.Range('AF165').Value = 'CL:0000625'
$ws.Range('AF166').Value = 'CL:0002397'
$ws.Range('AF167').Value = 'CL:0000236'
$ws.Range('AG167').Value = 'PR:000001002, PR:000001343'
$ws.Range('AF168').Value = 'CL:0002397'
$ws.Range('AF169').Value = 'CL:0000623'
$ws.Range('AG169').Value = 'PR:000001483, PR:000001024, PR:000001020 '
$ws.Range('AF170').Value = 'CL:0002397'
$ws.Range('AF171').Value = 'CL:0001054'
$ws.Range('AF172').Value = 'CL:0000786'
$ws.Range('AF173').Value = 'CL:0000576'
$ws.Range('AF174').Value = 'CL:0000084'
$ws.Range('AF175').Value = 'CL:0000786'
$ws.Range('AF176').Value = 'CL:0000798'
$ws.Range('AF177').Value = 'CL:0000625'
$ws.Range('AF178').Value = 'CL:0000625'
$ws.Range('AF179').Value = 'CL:0000798'
$ws.Range('AF180').Value = 'CL:0000798'
$ws.Range('AF181').Value = 'CL:0000798'
$ws.Range('AF182').Value = 'CL:0000798'
$ws.Range('AF183').Value = 'CL:0000798'
$ws.Range('AF184').Value = 'CL:0000625'
$ws.Range('AF185').Value = 'CL:0000904'
$ws.Range('AG185').Value = 'PR:000001017, PR:000001203'
$ws.Range('AF186').Value = 'CL:0000904'
$ws.Range('AG186').Value = 'PR:000001017, PR:000001203'
$ws.Range('AF187').Value = 'CL:0000904'
$ws.Range('AG187').Value = 'PR:000001017, PR:000001203'
$ws.Range('AF188').Value = 'CL:0000904'
$ws.Range('AG188').Value = 'PR:000001017, PR:000001203'
$ws.Range('AF189').Value = 'CL:0000904'
$ws.Range('AG189').Value = 'PR:000001017, PR:000001203'
$ws.Range('AF190').Value = 'CL:0000904'
$ws.Range('AG190').Value = 'PR:000001017, PR:000001203'
$ws.Range('AF191').Value = 'CL:0001041'
$ws.Range('AF192').Value = 'CL:0000980'
